# Translate the template spreadsheet's table (sheet) names from German to
# English. Renaming the sheets automatically keeps the `_xlnm._FilterDatabase`
# defined names (and their sheet-name prefix / quoting) in sync.
$wb = $excel.ActiveWorkbook

$overviewSheet = $wb.Worksheets.Item(1)
$monthSheet = $wb.Worksheets.Item(2)

$overviewSheet.Name = "Complete Overview"
$monthSheet.Name = "Month Overview"

# The author's last-saved view had the "Month Overview" tab active with
# C12 selected (previously the first sheet was active with A9 selected).
$monthSheet.Activate()
$monthSheet.Range("C12").Select()
